# ArrayPlan.xlsx - "Succesfully got board array printed"
#
# The board array in N4:U11 used to hold a "rank*10 + file" style numbering
# (21,22,23...28 / 31,32,...38 / ...). The sheet now prints a proper
# zero-based flattened board array (0..63) instead.
#
# The active selection on the sheet is also moved onto the board-array
# block that was just finished (N4:U11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Board array N4:U11: rank/file numbering -> flat 0..63 index ---------
$cols = @("N", "O", "P", "Q", "R", "S", "T", "U")
for ($r = 4; $r -le 11; $r++) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = ($r - 4) * $cols.Length + $i
        $ws.Range("$col$r").Value = $value
    }
}

# --- Move the active selection onto the board array that was finished ----
$ws.Range("N4:U11").Select() | Out-Null
